$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (coin names, links, percentage strings)
$plainValues = [ordered]@{
    "E2" = "  -1.71%  "
    "E3" = "  -1.86%  "
    "E4" = "  +0.39%  "
    "E5" = "  -1.26%  "
    "E6" = "  +0.34%  "
    "E7" = "  +1.55%  "
    "E8" = "  -1.40%  "
    "E9" = "  -2.76%  "
    "E10" = "  -2.37%  "
    "E11" = "  -2.42%  "
    "E12" = "  -1.32%  "
    "E13" = "  -3.44%  "
    "E14" = "  -0.98%  "
    "E15" = "  -0.82%  "
    "E16" = "  -4.26%  "
    "E18" = "  -2.01%  "
    "E19" = "  +0.32%  "
    "E20" = "  -2.75%  "
    "E21" = "  -1.69%  "
    "E22" = "  +0.52%  "
    "E23" = "  -2.73%  "
    "E24" = "  -1.04%  "
    "E25" = "  -4.00%  "
    "E26" = "  -1.05%  "
    "E27" = "  -9.23%  "
    "E28" = "  -1.94%  "
    "E29" = "  -2.62%  "
    "E30" = "  -3.71%  "
    "E31" = "  -0.02%  "
    "E32" = "  -1.96%  "
    "E33" = "  -2.12%  "
    "E34" = "  -3.11%  "
    "E35" = "  -0.67%  "
    "E36" = "  +0.32%  "
    "E37" = "  +1.12%  "
    "E38" = "  -2.54%  "
    "E39" = "  -1.54%  "
    "E40" = "  -3.53%  "
    "E41" = "  +0.23%  "
    "B42" = "RenderToken"
    "C42" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "E42" = "  +6.77%  "
    "B43" = "TheSandbox"
    "C43" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "E43" = "  -2.31%  "
    "E44" = "  -4.45%  "
    "E45" = "  -3.13%  "
    "E46" = "  -2.38%  "
    "B47" = "EnergySwap"
    "C47" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E47" = "  -3.79%  "
    "B48" = "PaxDollar"
    "C48" = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "E48" = "  +0.33%  "
    "E49" = "  -2.20%  "
    "E50" = "  -4.05%  "
    "E51" = "  -1.38%  "
}
foreach ($cellRef in $plainValues.Keys) {
    $ws.Range($cellRef).Value = $plainValues[$cellRef]
}

# Price column updates: force text storage so numeric-looking strings
# (e.g. "1.002", "26.660.62") are not reinterpreted as numbers, then
# restore the default "Normal" style so no stray number-format style lingers.
$priceValues = [ordered]@{
    "D2" = "26.660.62"
    "D3" = "1.786.45"
    "D5" = "308.03"
    "D6" = "1.002"
    "D7" = "0.4524"
    "D8" = "0.3687"
    "D9" = "0.07236"
    "D10" = "0.8513"
    "D12" = "1.794.81"
    "D13" = "6.505"
    "D14" = "5.280"
    "D15" = "0.07024"
    "D16" = "90.32"
    "D17" = "1.003"
    "D18" = "0.000008574"
    "D19" = "1.002"
    "D20" = "14.59"
    "D21" = "26.685.99"
    "D22" = "5.243"
    "D23" = "10.64"
    "D24" = "2.013.38"
    "D25" = "1.901"
    "D26" = "149.94"
    "D27" = "2.154"
    "D28" = "18.14"
    "D29" = "5.177"
    "D30" = "113.66"
    "D31" = "0.08815"
    "D32" = "0.7504"
    "D33" = "1.152"
    "D34" = "4.425"
    "D35" = "2.865"
    "D36" = "1.002"
    "D37" = "1.111"
    "D38" = "0.01934"
    "D39" = "0.05192"
    "D40" = "7.117"
    "D41" = "2.859"
    "D42" = "2.317"
    "D43" = "0.5184"
    "D45" = "8.426"
    "D46" = "0.4922"
    "D47" = "10.24"
    "D48" = "1.001"
    "D50" = "1.638"
    "D51" = "0.06280"
}
foreach ($cellRef in $priceValues.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceValues[$cellRef]
    $cell.Style = "Normal"
}
